$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: this workbook stores every cell as text, so force text format
# before assigning values to avoid Excel's automatic number coercion.

# Row 6: remove the empty placeholder cell at E6 entirely
$ws.Range("E6").ClearContents()

# Row 8: populate the row with the new set of values
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "2"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").ClearContents()

$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "5"

$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "no pero si"

$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "1"

# Row 10: update existing values
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "4"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "no pero si"
